$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "খাতা/পত্রের সংখ্যা" (number of scripts/papers) column (G)
# for several billing line items; the dependent formulas in column I
# (and the grand total in I32) recalculate automatically.
$ws.Range("G16").Value = 27
$ws.Range("G17").Value = 31
$ws.Range("G18").Value = 118
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1
